# Update the "想去人数" (interest count) figures in column F across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets to
# reflect the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 706
$ws1.Range("F5").Value = 2228
$ws1.Range("F6").Value = 1320
$ws1.Range("F10").Value = 2856
$ws1.Range("F13").Value = 1072
$ws1.Range("F16").Value = 916
$ws1.Range("F20").Value = 126
$ws1.Range("F26").Value = 4883
$ws1.Range("F28").Value = 151

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F26").Value = 570
$ws2.Range("F28").Value = 29
$ws2.Range("F36").Value = 713

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 706
$ws4.Range("F13").Value = 2228
$ws4.Range("F14").Value = 1320
$ws4.Range("F20").Value = 2856
$ws4.Range("F24").Value = 1072
$ws4.Range("F29").Value = 916
$ws4.Range("F30").Value = 916
$ws4.Range("F35").Value = 126
$ws4.Range("F45").Value = 4883
$ws4.Range("F48").Value = 151
$ws4.Range("F49").Value = 713
